$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '76.913.87'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.46%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.951.57'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +2.85%  '

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '201.32'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.48%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '596.54'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.550'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.62%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.196'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.35%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '2.937.43'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.32%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.449'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +14.49%  '

# Row 12
$ws.Range("E12").Value = '  +0.17%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.91'
$cell.NumberFormat = "General"
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '3.474.35'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +2.33%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '76.553.24'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.21%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '28.22'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.54%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.0000188'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.23%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '2.935.43'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.57%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '13.32'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +6.63%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '8.66'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -4.88%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '371.62'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -3.56%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.33'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +4.49%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '2.27'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -3.56%  '

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '72.34'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.66%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.35%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '3.080.64'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.40%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '4.23'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.70'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.92%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.0000107'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.04%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.21%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.19'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +5.77%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.38'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -2.90%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '494.37'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.97%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.83'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '166.60'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.16%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.113'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +23.31%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '20.17'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.44%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.394'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +14.16%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '19.76'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +1.28%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.110'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -7.51%  '

# Row 42
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '180.79'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.02%  '

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '4.90'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -3.86%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.65'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.16%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '40.14'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.09%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.18'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -4.75%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.590'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +1.66%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '3.88'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +3.55%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.07%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '22.71'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +4.70%  '
